$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 38891596
$ws.Cells.Item(70, 9).Value = 31252076
$ws.Cells.Item(70, 10).Value = 47622480
$ws.Cells.Item(70, 11).Value = 93756228
$ws.Cells.Item(70, 12).Value = 142867440
$ws.Cells.Item(70, 13).Value = -93755958
$ws.Cells.Item(70, 14).Value = -142867980

$ws.Cells.Item(73, 8).Value = 38891596
$ws.Cells.Item(73, 9).Value = 31252076
$ws.Cells.Item(73, 10).Value = 47622480
$ws.Cells.Item(73, 11).Value = 93756228
$ws.Cells.Item(73, 12).Value = 142867440
$ws.Cells.Item(73, 13).Value = -93755292
$ws.Cells.Item(73, 14).Value = -142869312

$ws.Cells.Item(112, 8).Value = 4899.4165
$ws.Cells.Item(112, 9).Value = 1404.3334
$ws.Cells.Item(112, 10).Value = 5132.4224
$ws.Cells.Item(112, 11).Value = 4213.0002
$ws.Cells.Item(112, 12).Value = 15397.2672
$ws.Cells.Item(112, 13).Value = -3105.0002
$ws.Cells.Item(112, 14).Value = -17613.2672

$ws.Cells.Item(129, 8).Value = 1791.5
$ws.Cells.Item(129, 9).Value = 1791.5
$ws.Cells.Item(129, 11).Value = 5374.5
$ws.Cells.Item(129, 13).Value = -374.5

$ws.Cells.Item(132, 8).Value = 2290.894
$ws.Cells.Item(132, 9).Value = 2189.258
$ws.Cells.Item(132, 10).Value = 3866.25
$ws.Cells.Item(132, 11).Value = 6567.773999999999
$ws.Cells.Item(132, 12).Value = 11598.75
$ws.Cells.Item(132, 13).Value = -4037.773999999999
$ws.Cells.Item(132, 14).Value = -16658.75

$ws.Cells.Item(138, 8).Value = 6546.8823
$ws.Cells.Item(138, 10).Value = 9204.454
$ws.Cells.Item(138, 12).Value = 27613.362
$ws.Cells.Item(138, 14).Value = -37893.362

$ws.Cells.Item(141, 8).Value = 2724.6667
$ws.Cells.Item(141, 9).Value = 2257.6
$ws.Cells.Item(141, 11).Value = 6772.799999999999
$ws.Cells.Item(141, 13).Value = -1592.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1714660
$ws.Cells.Item(32, 9).Value = 1840403.5
$ws.Cells.Item(32, 11).Value = 1840403.5
$ws.Cells.Item(32, 13).Value = -1840116.5

$ws.Cells.Item(61, 8).Value = 10681.412
$ws.Cells.Item(61, 9).Value = 2765
$ws.Cells.Item(61, 10).Value = 14999.454
$ws.Cells.Item(61, 11).Value = 2765
$ws.Cells.Item(61, 12).Value = 14999.454
$ws.Cells.Item(61, 13).Value = -2553
$ws.Cells.Item(61, 14).Value = -15423.454

$ws.Cells.Item(110, 8).Value = 27778790
$ws.Cells.Item(110, 9).Value = 1053.625
$ws.Cells.Item(110, 10).Value = 83334264
$ws.Cells.Item(110, 11).Value = 1053.625
$ws.Cells.Item(110, 12).Value = 83334264
$ws.Cells.Item(110, 13).Value = 991.375
$ws.Cells.Item(110, 14).Value = -83338354

$ws.Cells.Item(122, 8).Value = 14736.823
$ws.Cells.Item(122, 9).Value = 17425.154
$ws.Cells.Item(122, 11).Value = 52275.462
$ws.Cells.Item(122, 13).Value = -49825.462

$ws.Cells.Item(132, 8).Value = 6408.2856
$ws.Cells.Item(132, 9).Value = 3060.889
$ws.Cells.Item(132, 11).Value = 9182.667000000001
$ws.Cells.Item(132, 13).Value = -6652.667000000001

$ws.Cells.Item(136, 8).Value = 10681.412
$ws.Cells.Item(136, 9).Value = 2765
$ws.Cells.Item(136, 10).Value = 14999.454
$ws.Cells.Item(136, 11).Value = 8295
$ws.Cells.Item(136, 12).Value = 44998.362
$ws.Cells.Item(136, 13).Value = -5745
$ws.Cells.Item(136, 14).Value = -50098.362

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2951.88
$ws.Cells.Item(105, 9).Value = 2205.75
$ws.Cells.Item(105, 11).Value = 2205.75
$ws.Cells.Item(105, 13).Value = -458.75

$ws.Cells.Item(134, 8).Value = 5394.959
$ws.Cells.Item(134, 9).Value = 2634.7097
$ws.Cells.Item(134, 10).Value = 10148.723
$ws.Cells.Item(134, 11).Value = 7904.1291
$ws.Cells.Item(134, 12).Value = 30446.169
$ws.Cells.Item(134, 13).Value = -5369.1291
$ws.Cells.Item(134, 14).Value = -35516.169

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6574.5645
$ws.Cells.Item(31, 9).Value = 2895.7144
$ws.Cells.Item(31, 11).Value = 2895.7144
$ws.Cells.Item(31, 13).Value = -2600.7144

$ws.Cells.Item(34, 8).Value = 6574.5645
$ws.Cells.Item(34, 9).Value = 2895.7144
$ws.Cells.Item(34, 11).Value = 2895.7144
$ws.Cells.Item(34, 13).Value = -2693.7144

$ws.Cells.Item(59, 8).Value = 97137.57000000001
$ws.Cells.Item(59, 10).Value = 97137.57000000001
$ws.Cells.Item(59, 12).Value = 97137.57000000001
$ws.Cells.Item(59, 14).Value = -99427.57000000001

$ws.Cells.Item(132, 8).Value = 4663.967
$ws.Cells.Item(132, 9).Value = 1780.2632
$ws.Cells.Item(132, 11).Value = 5340.7896
$ws.Cells.Item(132, 13).Value = -2810.7896

$ws.Cells.Item(134, 8).Value = 8849.77
$ws.Cells.Item(134, 9).Value = 9636.647000000001
$ws.Cells.Item(134, 10).Value = 8241.727999999999
$ws.Cells.Item(134, 11).Value = 28909.941
$ws.Cells.Item(134, 12).Value = 24725.184
$ws.Cells.Item(134, 13).Value = -26374.941
$ws.Cells.Item(134, 14).Value = -29795.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 5004537.5
$ws.Cells.Item(5, 10).Value = 6299.8
$ws.Cells.Item(5, 12).Value = 18899.4
$ws.Cells.Item(5, 14).Value = -19123.4

$ws.Cells.Item(12, 8).Value = 2273471.5
$ws.Cells.Item(12, 10).Value = 3572036.5
$ws.Cells.Item(12, 12).Value = 10716109.5
$ws.Cells.Item(12, 14).Value = -10716455.5

$ws.Cells.Item(14, 8).Value = 15152103
$ws.Cells.Item(14, 9).Value = 15152103
$ws.Cells.Item(14, 11).Value = 45456309
$ws.Cells.Item(14, 13).Value = -45456136

$ws.Cells.Item(92, 8).Value = 19232016
$ws.Cells.Item(92, 10).Value = 19232016
$ws.Cells.Item(92, 12).Value = 57696048
$ws.Cells.Item(92, 14).Value = -57698544

$ws.Cells.Item(98, 8).Value = 472.57144
$ws.Cells.Item(98, 10).Value = 666.6667
$ws.Cells.Item(98, 12).Value = 2000.0001
$ws.Cells.Item(98, 14).Value = -4996.0001

$ws.Cells.Item(107, 8).Value = 15385445
$ws.Cells.Item(107, 10).Value = 20000868
$ws.Cells.Item(107, 12).Value = 60002604
$ws.Cells.Item(107, 14).Value = -60006444

$ws.Cells.Item(121, 8).Value = 20000578
$ws.Cells.Item(121, 9).Value = 20000278
$ws.Cells.Item(121, 10).Value = 20000878
$ws.Cells.Item(121, 11).Value = 60000834
$ws.Cells.Item(121, 12).Value = 60002634
$ws.Cells.Item(121, 13).Value = -59999524
$ws.Cells.Item(121, 14).Value = -60005254

$ws.Cells.Item(132, 8).Value = 10290.056
$ws.Cells.Item(132, 10).Value = 12318.909
$ws.Cells.Item(132, 12).Value = 110870.181
$ws.Cells.Item(132, 14).Value = -115930.181

$ws.Cells.Item(135, 8).Value = 5004537.5
$ws.Cells.Item(135, 10).Value = 6299.8
$ws.Cells.Item(135, 12).Value = 56698.2
$ws.Cells.Item(135, 14).Value = -61768.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 95433.55
$ws.Cells.Item(122, 9).Value = 202714.4
$ws.Cells.Item(122, 11).Value = 608143.2
$ws.Cells.Item(122, 13).Value = -605693.2

$ws.Cells.Item(126, 8).Value = 2663.2856
$ws.Cells.Item(126, 9).Value = 2918.6
$ws.Cells.Item(126, 10).Value = 2521.4443
$ws.Cells.Item(126, 11).Value = 8755.799999999999
$ws.Cells.Item(126, 12).Value = 7564.3329
$ws.Cells.Item(126, 13).Value = -6285.799999999999
$ws.Cells.Item(126, 14).Value = -12504.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 14561207
$ws.Cells.Item(46, 9).Value = 8621365
$ws.Cells.Item(46, 11).Value = 8621365
$ws.Cells.Item(46, 13).Value = -8621177

$ws.Cells.Item(61, 8).Value = 6330.3
$ws.Cells.Item(61, 9).Value = 4251.5
$ws.Cells.Item(61, 10).Value = 6850
$ws.Cells.Item(61, 11).Value = 4251.5
$ws.Cells.Item(61, 12).Value = 6850
$ws.Cells.Item(61, 13).Value = -4049.5
$ws.Cells.Item(61, 14).Value = -7254

$ws.Cells.Item(113, 8).Value = 6330.3
$ws.Cells.Item(113, 9).Value = 4251.5
$ws.Cells.Item(113, 11).Value = 4251.5
$ws.Cells.Item(113, 13).Value = -2081.5
$ws.Cells.Item(113, 14).Value = -11190

$ws.Cells.Item(127, 8).Value = 61683.168
$ws.Cells.Item(127, 10).Value = 61683.168
$ws.Cells.Item(127, 12).Value = 61683.168
$ws.Cells.Item(127, 14).Value = -71603.16800000001
